$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.972.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = "'1.676.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'215.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("E6").Value = '  +1.30%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").Value = "'20.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").Value = "'0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = "'1.912.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = "'1.671.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = "'0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = "'65.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = "'26.984.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = "'236.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = "'8.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.16%  '
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").Value = "'9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = "'145.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").Value = "'7.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").Value = "'1.485.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("E35").Value = '  +4.26%  '
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = "'0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").Value = "'0.0174"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.83%  '
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("E41").Value = '  +1.23%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("D44").Value = "'67.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("D45").Value = "'1.818.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = "'0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = "'90.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("E51").Value = '  +0.38%  '
